$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Veda")

# --- Scenario-category lookup table (A7:C61) -------------------------------
# The first block (rows 7-11) used to hold five distinct literal category
# labels ("C1","C2","C3","C4","C7"); the user instead filled the formula
# from C7 down through C11, so the whole first block now resolves to "C1"
# (matching the pre-existing =<row above> pattern already used by the rows
# further down the table).
$ws.Range("C7").AutoFill($ws.Range("C7:C11"), 0)  # 0 = xlFillDefault (fill down)

# --- fuel_prices block header / formula fix (rows 42-46) -------------------
# F42 header label changed from "process" to "pset_co".
$ws.Range("F42").Value = "pset_co"

# F43 used to be hard-coded to the literal text "fuel_supply"; it should
# instead mirror the commodity name in column E. Enter the formula in F43
# and fill it down through F46.
$ws.Range("F43").Formula = "=E43"
$ws.Range("F43").AutoFill($ws.Range("F43:F46"), 0)

# --- Final view state --------------------------------------------------
$ws.Range("F47").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

$excel.CalculateFull()
